$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 for the "Verify Resend OTP button" test case.
# This shifts the former rows 7-10 down to 8-11.
$ws.Rows.Item(7).Insert()

# Rewrite rows 2-11 with updated sequence numbers, text and refreshed execution timestamps.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "User Login with Valid Credentials"
$ws.Range("C2").Value = "PASSED"
$ws.Range("D2").Value = "02/04/2025 01:07:39 PM"
$ws.Range("E2").Value = "Test executed successfully."

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Verify empty email state"
$ws.Range("C3").Value = "PASSED"
$ws.Range("D3").Value = "02/04/2025 01:07:42 PM"
$ws.Range("E3").Value = "Test executed successfully."

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Verify Login with Invalid Email"
$ws.Range("C4").Value = "PASSED"
$ws.Range("D4").Value = "02/04/2025 01:07:43 PM"
$ws.Range("E4").Value = "Test executed successfully."

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Verify Login with Not Registred Email"
$ws.Range("C5").Value = "PASSED"
$ws.Range("D5").Value = "02/04/2025 01:07:44 PM"
$ws.Range("E5").Value = "Test executed successfully."

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Verify Login with Invalid OTP"
$ws.Range("C6").Value = "PASSED"
$ws.Range("D6").Value = "02/04/2025 01:07:46 PM"
$ws.Range("E6").Value = "Test executed successfully."

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Verify Resend OTP button"
$ws.Range("C7").Value = "PASSED"
$ws.Range("D7").Value = "02/04/2025 01:08:48 PM"
$ws.Range("E7").Value = "Test executed successfully."

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Verify account block after attempting wrong OTP for 5 times"
$ws.Range("C8").Value = "PASSED"
$ws.Range("D8").Value = "02/04/2025 01:09:00 PM"
$ws.Range("E8").Value = "Test executed successfully."

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Verify Go To Sign In page Navigation"
$ws.Range("C9").Value = "PASSED"
$ws.Range("D9").Value = "02/04/2025 01:09:02 PM"
$ws.Range("E9").Value = "Test executed successfully."

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Verify that navigation and getOTP blocked for blocked account"
$ws.Range("C10").Value = "PASSED"
$ws.Range("D10").Value = "02/04/2025 01:09:03 PM"
$ws.Range("E10").Value = "Test executed successfully."

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Verify Home Page Loads Successfully"
$ws.Range("C11").Value = "PASSED"
$ws.Range("D11").Value = "02/04/2025 01:09:03 PM"
$ws.Range("E11").Value = "Test executed successfully."
